$d = $word.ActiveDocument

# --- Insert the two new paragraphs after the first empty NormalWeb paragraph (#3) ---
$d.Paragraphs(3).Range.InsertParagraphAfter()
$d.Paragraphs(4).Range.InsertAfter("Getting locs is not just a change of hairstyle but also a change of lifestyle. During a traditional maintenance session, it will take a loctician 45 mins to 1 hour and a half to palm roll a full set of traditional locs. Sisterlocs are microlocs that are interlocked using a special tool. When sisterlocs are installed, depending on your hair length, your hair may be interlocked down the length of your hair, but definitely at the root. This allows you to have the locs without having to go through the different phases of locking your hair. Therefore, they are perfect for people who would find it challenging to grow traditional locks and prefer quick results.")

$d.Paragraphs(4).Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.InsertAfter("Washing and conditioning our natural hair can be a task that no one wants to worry about. It’s understandable if you are maybe unsure about the best way to wash your hair, or you wonder what the best shampoo for natural hair is, or the best hair wash routine. Kandase uses the right techniques that can significantly improve your hair structure to achieve the maximum shine, the largest possible volume and the long-lasting health of your hair.")

# --- Insert one new paragraph after the next empty NormalWeb paragraph (originally #4, now #6) ---
$d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.InsertAfter("Protective styling is one of the best and most efficient ways to achieve longer hair, strengthen your hair, and keep it healthy. The most popular go-to styles for naturalistas when they want a protective style are braids and twists. ")

# --- Insert five new empty NormalWeb paragraphs after the last empty paragraph (originally #5, now #8) ---
$d.Paragraphs(8).Range.InsertParagraphAfter()
$d.Paragraphs(9).Range.InsertParagraphAfter()
$d.Paragraphs(10).Range.InsertParagraphAfter()
$d.Paragraphs(11).Range.InsertParagraphAfter()
$d.Paragraphs(12).Range.InsertParagraphAfter()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
